$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column A (date strings) to Text format so Excel
# doesn't auto-convert the dd-mm-yyyy-looking strings into date serial
# numbers when we assign them.
$ws.Range("A3:A21").NumberFormat = "@"

# Update date strings from slash-format to dash-format
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("A6").Value = "08-08-2022"
$ws.Range("A7").Value = "11-08-2022"
$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

# Restore the original (default) style now that the text values are safely set.
$ws.Range("A3:A21").Style = "Normal"

# Update the attendance counters that changed in the recount
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0

$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("H9").Value = 0

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0
